# Add team record (Wins / Losses / Ties) columns to the roster sheet.
# New columns AD:AF are appended after the existing data (A:AC), with a
# bold/centered/bordered header row matching the look of the existing
# header, and a constant record (66-96-0) filled in for every player row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the formatting already used for the other header cells
# (bold font, centered/top aligned, thin border all around).
$hdrRange = $ws.Range("AD1:AF1")
$hdrRange.Font.Bold = $true
$hdrRange.HorizontalAlignment = -4108   # xlCenter
$hdrRange.VerticalAlignment = -4160     # xlTop
$hdrRange.Borders.LineStyle = 1         # xlContinuous

# --- Data rows (2 through 47): team record for every player ---
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 47) { $lastRow = 47 }

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 66  # AD: Wins
    $ws.Cells.Item($r, 31).Value = 96  # AE: Losses
    $ws.Cells.Item($r, 32).Value = 0   # AF: Ties
}
